$wb = $excel.ActiveWorkbook

$schemaSheet = $wb.Worksheets.Item("!!_Schema")
$compoundSheet = $wb.Worksheets.Item("!!Compound")
$modelSheet = $wb.Worksheets.Item("!!Model")
$reactionSheet = $wb.Worksheets.Item("!!Reaction")

$schemaSheet.Unprotect()
$compoundSheet.Unprotect()
$modelSheet.Unprotect()
$reactionSheet.Unprotect()

$schemaSheet.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8' date='2020-03-09 12:59:51'"
$schemaSheet.Range("A2").Value = "!!ObjTables type='Schema' description='Table/model and column/attribute definitions' date='2020-03-09 12:59:51' objTablesVersion='0.0.8'"

$compoundSheet.Range("A1").Value = "!!ObjTables type='Data' id='Compound' description='Compound' name='Compound' date='2020-03-09 12:59:51' objTablesVersion='0.0.8'"

$modelSheet.Range("A1").Value = "!!ObjTables type='Data' id='Model' description='Model' name='Model' date='2020-03-09 12:59:51' objTablesVersion='0.0.8'"

$reactionSheet.Range("A1").Value = "!!ObjTables type='Data' id='Reaction' description='Reaction' name='Reaction' date='2020-03-09 12:59:51' objTablesVersion='0.0.8'"

$schemaSheet.Protect($null, $true, $true, $true, $false, $true, $true, $true, $true, $false, $true, $true, $false, $true, $true, $true)
$compoundSheet.Protect($null, $true, $true, $true, $false, $true, $true, $true, $true, $false, $true, $true, $false, $true, $true, $true)
$modelSheet.Protect($null, $true, $true, $true, $false, $true, $true, $true, $true, $false, $true, $true, $false, $true, $true, $true)
$reactionSheet.Protect($null, $true, $true, $true, $false, $true, $true, $true, $true, $false, $true, $true, $false, $true, $true, $true)
